$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '35.412.81'
$ws.Range('E2').Value = '  +1.85%  '
$ws.Range('D3').Value = '1.886.32'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '246.31'
$ws.Range('E5').Value = '  -0.04%  '
$ws.Range('B6').Value = 'XRP'
$ws.Range('C6').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.692'
$ws.Range('E6').Value = '  +0.94%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '43.30'
$ws.Range('E8').Value = '  +3.15%  '
$ws.Range('E9').Value = '  +2.69%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '54.86'
$ws.Range('E10').Value = '  +7.76%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0744'
$ws.Range('E11').Value = '  +1.36%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.0988'
$ws.Range('E12').Value = '  +1.91%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '13.84'
$ws.Range('E13').Value = '  +8.76%  '
$ws.Range('D14').Value = '2.156.44'
$ws.Range('E14').Value = '  +0.69%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.771'
$ws.Range('E15').Value = '  +8.28%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.02'
$ws.Range('E16').Value = '  +2.81%  '
$ws.Range('D17').Value = '1.903.51'
$ws.Range('E17').Value = '  +1.01%  '
$ws.Range('D18').Value = '35.368.42'
$ws.Range('E18').Value = '  +1.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '73.52'
$ws.Range('E19').Value = '  +1.25%  '
$ws.Range('D20').Value = '0.0₃0826'
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '244.62'
$ws.Range('E21').Value = '  +0.69%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '12.83'
$ws.Range('E22').Value = '  +1.43%  '
$ws.Range('E23').Value = '  +4.66%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.63'
$ws.Range('E24').Value = '  +8.22%  '
$ws.Range('E25').Value = '  +0.05%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.16'
$ws.Range('E26').Value = '  -3.31%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '165.62'
$ws.Range('E27').Value = '  +0.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '8.62'
$ws.Range('E28').Value = '  +3.26%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '18.29'
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.128'
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0596'
$ws.Range('E31').Value = '  +3.78%  '
$ws.Range('E32').Value = '  +0.62%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.88'
$ws.Range('E33').Value = '  +21.75%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.17'
$ws.Range('E34').Value = '  +0.86%  '
$ws.Range('E35').Value = '  -0.06%  '
$ws.Range('E36').Value = '  -13.74%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.857'
$ws.Range('E37').Value = '  +3.80%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.95'
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0720'
$ws.Range('E39').Value = '  +9.45%  '
$ws.Range('E40').Value = '  +5.36%  '
$ws.Range('B41').Value = 'InjectiveProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '17.27'
$ws.Range('E41').Value = '  +3.03%  '
$ws.Range('B42').Value = 'Aave'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '97.77'
$ws.Range('E42').Value = '  +0.41%  '
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '13.80'
$ws.Range('E44').Value = '  +14.46%  '
$ws.Range('D45').Value = '1.324.30'
$ws.Range('E45').Value = '  +3.43%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.39'
$ws.Range('E46').Value = '  +3.30%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.0810'
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('E49').Value = '  +0.50%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '6.28'
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('D51').Value = '2.056.02'
$ws.Range('E51').Value = '  +0.31%  '
